$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.090022772380280003,
    -0.0059999999481554767,
    0.018620277782507522,
    -0.0079999999253299592,
    -0.0029999999668222088,
    -0.0019999999743127717,
    -0.0099999999073414614,
    -0.0099999999080830904,
    -0.0019999999792332801,
    -0.0019999999825888182,
    -0.002999999974295342,
    -0.0034999999710674246,
    -0.0034999999776337276,
    0.035120838623606154,
    0.0000070942426608766596,
    -0.0019999999922495348,
    -0.0019999999915665256,
    -0.003999999974335644,
    -0.00399999996494671,
    -0.0039999999621773696,
    -0.0039999999617243986,
    -0.00399999996136291,
    -0.0049999999493195446,
    -0.01999999981626388,
    -0.019999999813501201,
    -0.0024999999698422926,
    -0.0024999999695545227,
    -0.0019999999727238205,
    -0.0069999999297269966,
    0.019976916534831801,
    -0.0069999999267906787,
    -0.0099999999008861806,
    -0.0039999999518194329
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Target stored width is 15.5703125 characters; this engine quantizes
# ColumnWidth to 1/6-character steps, so 14.6666... is the input that
# reliably lands on the nearest achievable stored width (15.5).
$ws.Columns.Item(1).ColumnWidth = 14.666666666666666
